# HRMS DATA ADDED SUCCESSFULLY
# Adds a new "LeaveConfiguration" worksheet (after the existing "BranchMaster"
# sheet) containing Business Unit / Week / WeekEnd configuration data, and
# makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet (BranchMaster)
# so it becomes sheet #6, at the end of the tab strip.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "LeaveConfiguration"

# Header row
$newSheet.Range("A1").Value = "Bunit"
$newSheet.Range("B1").Value = "WEEK"
$newSheet.Range("C1").Value = "WeekEND"

# Data row
$newSheet.Range("A2").Value = "BU1-Test"
$newSheet.Range("B2").Value = "Saturday"
$newSheet.Range("C2").Value = "Sunday"

# Column widths matching the source sheet's layout
$newSheet.Columns.Item(1).ColumnWidth = 16.833333333333332
$newSheet.Columns.Item(2).ColumnWidth = 22.666666666666668
$newSheet.Columns.Item(3).ColumnWidth = 20.166666666666668

# Centered alignment for the populated data range
$dataRange = $newSheet.Range("A1:C2")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108

# Portrait page orientation for printing
$newSheet.PageSetup.Orientation = 1

# Make this the active sheet/selection, matching the authored workbook state
[void]$newSheet.Range("C11").Select()
